# Insert two new weekly price records for "Camote" (Zapallo) at the top of
# the data block that starts at row 205, shifting all existing records
# (old rows 205-331) down by two rows (to new rows 207-333). The new rows
# carry the latest "guarda" quality prices dated 44438.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 205 onward down by two rows.
$ws.Rows("205:206").Insert()

# New row 205: "1a (guarda)"
$ws.Range("A205").Value = 8
$ws.Range("B205").Value = "Terminal La Palmera de La Serena"
$ws.Range("C205").Value = "Coquimbo"
$ws.Range("D205").Value = 44438
$ws.Range("E205").Value = 4
$ws.Range("F205").Value = 100112045
$ws.Range("G205").Value = "Zapallo"
$ws.Range("H205").Value = "Camote"
$ws.Range("I205").Value = "1a (guarda)"
$ws.Range("J205").Value = 800
$ws.Range("K205").Value = 1050
$ws.Range("L205").Value = 1100
$ws.Range("M205").Value = 1075
$ws.Range("N205").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O205").Value = "Región de O'Higgins"
$ws.Range("P205").Value = 1075
$ws.Range("Q205").Value = 1
$ws.Range("R205").Value = "Hortaliza"

# New row 206: "2a (guarda)"
$ws.Range("A206").Value = 8
$ws.Range("B206").Value = "Terminal La Palmera de La Serena"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 44438
$ws.Range("E206").Value = 4
$ws.Range("F206").Value = 100112045
$ws.Range("G206").Value = "Zapallo"
$ws.Range("H206").Value = "Camote"
$ws.Range("I206").Value = "2a (guarda)"
$ws.Range("J206").Value = 560
$ws.Range("K206").Value = 900
$ws.Range("L206").Value = 1000
$ws.Range("M206").Value = 950
$ws.Range("N206").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O206").Value = "Región de O'Higgins"
$ws.Range("P206").Value = 950
$ws.Range("Q206").Value = 1
$ws.Range("R206").Value = "Hortaliza"
